$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.382.28"
$ws.Range("E2").Value = "  +0.39%  "
$ws.Range("D3").Value = "2.576.14"
$ws.Range("E3").Value = "  +0.11%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'584.59"
$ws.Range("E5").Value = "  +0.56%  "
$ws.Range("D6").Value = "'173.67"
$ws.Range("E6").Value = "  +1.58%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'0.520"
$ws.Range("E8").Value = "  +1.52%  "
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Value = "'0.168"
$ws.Range("E9").Value = "  +0.86%  "
$ws.Range("B10").Value = "LidoStakedEther"
$ws.Range("C10").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D10").Value = "2.574.74"
$ws.Range("E10").Value = "  -0.15%  "
$ws.Range("E11").Value = "  -0.33%  "
$ws.Range("D12").Value = "'0.359"
$ws.Range("E12").Value = "  +2.14%  "
$ws.Range("E13").Value = "  +1.33%  "
$ws.Range("D14").Value = "3.049.20"
$ws.Range("E14").Value = "  -0.19%  "
$ws.Range("D15").Value = "71.253.47"
$ws.Range("E15").Value = "  +0.29%  "
$ws.Range("D16").Value = "'0.0000181"
$ws.Range("E16").Value = "  -1.66%  "
$ws.Range("D17").Value = "'25.63"
$ws.Range("E17").Value = "  +1.54%  "
$ws.Range("D18").Value = "2.581.14"
$ws.Range("E18").Value = "  -1.23%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "'11.68"
$ws.Range("E19").Value = "  -1.64%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "'7.98"
$ws.Range("E20").Value = "  +3.63%  "
$ws.Range("D21").Value = "'359.27"
$ws.Range("E21").Value = "  -1.62%  "
$ws.Range("D22").Value = "'3.98"
$ws.Range("E22").Value = "  -0.65%  "
$ws.Range("E23").Value = "  +4.97%  "
$ws.Range("D24").Value = "'1.00"
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("D25").Value = "'70.63"
$ws.Range("E25").Value = "  +0.49%  "
$ws.Range("D26").Value = "'4.15"
$ws.Range("E26").Value = "  -0.33%  "
$ws.Range("D27").Value = "'9.23"
$ws.Range("E27").Value = "  -0.17%  "
$ws.Range("D28").Value = "2.686.52"
$ws.Range("E28").Value = "  -2.50%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.13%  "
$ws.Range("D30").Value = "0.0₃0932"
$ws.Range("E30").Value = "  +0.23%  "
$ws.Range("D31").Value = "'8.02"
$ws.Range("E31").Value = "  +2.60%  "
$ws.Range("D32").Value = "'478.49"
$ws.Range("E32").Value = "  -1.31%  "
$ws.Range("E33").Value = "  -2.04%  "
$ws.Range("D34").Value = "'1.79"
$ws.Range("E34").Value = "  +0.22%  "
$ws.Range("E35").Value = "  -0.21%  "
$ws.Range("D36").Value = "'0.119"
$ws.Range("E36").Value = "  +3.29%  "
$ws.Range("D37").Value = "'158.25"
$ws.Range("E37").Value = "  +0.27%  "
$ws.Range("B38").Value = "EthereumClassic"
$ws.Range("C38").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D38").Value = "'18.89"
$ws.Range("E38").Value = "  -0.05%  "
$ws.Range("B39").Value = "WhiteBITCoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D39").Value = "'19.12"
$ws.Range("E39").Value = "  +1.28%  "
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("D41").Value = "'4.94"
$ws.Range("E41").Value = "  +3.27%  "
$ws.Range("E42").Value = "  +1.32%  "
$ws.Range("E43").Value = "  -3.42%  "
$ws.Range("D44").Value = "'2.41"
$ws.Range("E44").Value = "  -3.28%  "
$ws.Range("D45").Value = "'1.19"
$ws.Range("E45").Value = "  -11.10%  "
$ws.Range("D46").Value = "'38.70"
$ws.Range("E46").Value = "  +0.09%  "
$ws.Range("D47").Value = "'146.59"
$ws.Range("E47").Value = "  -0.54%  "
$ws.Range("D48").Value = "'0.543"
$ws.Range("E48").Value = "  +2.24%  "
$ws.Range("D49").Value = "'3.58"
$ws.Range("E49").Value = "  -0.36%  "
$ws.Range("E50").Value = "  -0.23%  "
$ws.Range("D51").Value = "'0.0743"
$ws.Range("E51").Value = "  +0.87%  "
